$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Models")

# Add a new row of data (row 3): Desktop / hardware / 1300 / 20 / Red
$ws.Range("A3").Value = "Desktop"
$ws.Range("B3").Value = "hardware"
$ws.Range("C3").Value = 1300
$ws.Range("D3").Value = 20
$ws.Range("E3").Value = "Red"

# Change the existing "Blue" description (row 2, col E) to "Pink"
$ws.Range("E2").Value = "Pink"

# Update the selected cell to match the saved view state
$ws.Range("E8").Select()
